$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1200
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1025
$ws.Range("N40").ClearContents()

$ws.Range("H41").Value = 475.3
$ws.Range("I41").Value = 406.375
$ws.Range("J41").Value = 751
$ws.Range("K41").Value = 406.375
$ws.Range("L41").Value = 751
$ws.Range("M41").Value = 33.625
$ws.Range("N41").Value = -1631

$ws.Range("H51").Value = 10865.15
$ws.Range("I51").Value = 4660
$ws.Range("J51").Value = 11960.177
$ws.Range("K51").Value = 4660
$ws.Range("L51").Value = 11960.177
$ws.Range("M51").Value = -4176
$ws.Range("N51").Value = -12928.177

$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 5400
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 16200
$ws.Range("M70").Value = -11730
$ws.Range("N70").Value = -16740

$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 5400
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 16200
$ws.Range("M73").Value = -11064
$ws.Range("N73").Value = -18072

$ws.Range("H92").Value = 965.1539
$ws.Range("I92").Value = 685.63635
$ws.Range("K92").Value = 685.63635
$ws.Range("M92").Value = 562.36365

$ws.Range("H96").Value = 259.875
$ws.Range("I96").Value = 230.64285
$ws.Range("K96").Value = 691.9285500000001
$ws.Range("M96").Value = 681.0714499999999

$ws.Range("H99").Value = 892.06665
$ws.Range("J99").Value = 1503
$ws.Range("L99").Value = 4509
$ws.Range("N99").Value = -7505

$ws.Range("H100").Value = 2357.5715
$ws.Range("J100").Value = 5499.5
$ws.Range("L100").Value = 5499.5
$ws.Range("N100").Value = -6581.5

$ws.Range("H123").Value = 68962.5
$ws.Range("J123").Value = 68962.5
$ws.Range("L123").Value = 68962.5
$ws.Range("N123").Value = -78762.5

$ws.Range("H132").Value = 1291918
$ws.Range("I132").Value = 2388.4
$ws.Range("J132").Value = 6127654
$ws.Range("K132").Value = 7165.200000000001
$ws.Range("L132").Value = 18382962
$ws.Range("M132").Value = -4635.200000000001
$ws.Range("N132").Value = -18388022

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1898.9333
$ws.Range("I2").Value = 1669.1
$ws.Range("J2").Value = 2358.6
$ws.Range("K2").Value = 1669.1
$ws.Range("L2").Value = 2358.6
$ws.Range("M2").Value = -1556.1
$ws.Range("N2").Value = -2584.6

$ws.Range("H32").Value = 904.99
$ws.Range("I32").Value = 790.1325000000001
$ws.Range("J32").Value = 1465.7646
$ws.Range("K32").Value = 790.1325000000001
$ws.Range("L32").Value = 1465.7646
$ws.Range("M32").Value = -503.1325000000001
$ws.Range("N32").Value = -2039.7646

$ws.Range("H97").Value = 2232975.8
$ws.Range("I97").Value = 3125796
$ws.Range("K97").Value = 3125796
$ws.Range("M97").Value = -3125300

$ws.Range("H102").Value = 8410995
$ws.Range("I102").Value = 10212708
$ws.Range("K102").Value = 10212708
$ws.Range("M102").Value = -10211086

$ws.Range("H116").Value = 1898.9333
$ws.Range("I116").Value = 1669.1
$ws.Range("J116").Value = 2358.6
$ws.Range("K116").Value = 1669.1
$ws.Range("L116").Value = 2358.6
$ws.Range("M116").Value = 624.9000000000001
$ws.Range("N116").Value = -6946.6

$ws.Range("H122").Value = 4117269.8
$ws.Range("I122").Value = 2160.4285
$ws.Range("J122").Value = 18520152
$ws.Range("K122").Value = 6481.2855
$ws.Range("L122").Value = 55560456
$ws.Range("M122").Value = -4031.2855
$ws.Range("N122").Value = -55565356

$ws.Range("H132").Value = 96268
$ws.Range("I132").Value = 67624.47
$ws.Range("J132").Value = 157647
$ws.Range("K132").Value = 202873.41
$ws.Range("L132").Value = 472941
$ws.Range("M132").Value = -200343.41
$ws.Range("N132").Value = -478001

$ws.Range("H133").Value = 33195.75
$ws.Range("J133").Value = 33195.75
$ws.Range("L133").Value = 33195.75
$ws.Range("N133").Value = -38255.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1898.9333
$ws.Range("I3").Value = 1669.1
$ws.Range("J3").Value = 2358.6
$ws.Range("K3").Value = 1669.1
$ws.Range("L3").Value = 2358.6
$ws.Range("M3").Value = -1555.1
$ws.Range("N3").Value = -2586.6

$ws.Range("H86").Value = 12109.167
$ws.Range("I86").Value = 26931
$ws.Range("J86").Value = 3216.0667
$ws.Range("K86").Value = 26931
$ws.Range("L86").Value = 3216.0667
$ws.Range("M86").Value = -25808
$ws.Range("N86").Value = -5462.066699999999

$ws.Range("H89").Value = 12109.167
$ws.Range("I89").Value = 26931
$ws.Range("J89").Value = 3216.0667
$ws.Range("K89").Value = 134655
$ws.Range("L89").Value = 16080.3335
$ws.Range("M89").Value = -129039
$ws.Range("N89").Value = -27312.3335

$ws.Range("H94").Value = 740.9474
$ws.Range("I94").Value = 680.2727
$ws.Range("J94").Value = 824.375
$ws.Range("K94").Value = 680.2727
$ws.Range("L94").Value = 824.375
$ws.Range("M94").Value = -229.2727
$ws.Range("N94").Value = -1726.375

$ws.Range("H99").Value = 1275.238
$ws.Range("I99").Value = 1210.5
$ws.Range("J99").Value = 1663.6666
$ws.Range("K99").Value = 1210.5
$ws.Range("L99").Value = 1663.6666
$ws.Range("M99").Value = 287.5
$ws.Range("N99").Value = -4659.6666

$ws.Range("H105").Value = 45457500
$ws.Range("I105").Value = 71431510
$ws.Range("J105").Value = 2975
$ws.Range("K105").Value = 71431510
$ws.Range("L105").Value = 2975
$ws.Range("M105").Value = -71429763
$ws.Range("N105").Value = -6469

$ws.Range("H123").Value = 23749
$ws.Range("J123").Value = 23749
$ws.Range("L123").Value = 23749
$ws.Range("N123").Value = -33549

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1796.2924
$ws.Range("I31").Value = 1216.4894
$ws.Range("J31").Value = 3310.2222
$ws.Range("K31").Value = 1216.4894
$ws.Range("L31").Value = 3310.2222
$ws.Range("M31").Value = -921.4893999999999
$ws.Range("N31").Value = -3900.2222

$ws.Range("H34").Value = 1796.2924
$ws.Range("I34").Value = 1216.4894
$ws.Range("J34").Value = 3310.2222
$ws.Range("K34").Value = 1216.4894
$ws.Range("L34").Value = 3310.2222
$ws.Range("M34").Value = -1014.4894
$ws.Range("N34").Value = -3714.2222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -10340

$ws.Range("H131").Value = 1005.89655
$ws.Range("J131").Value = 1132.2113
$ws.Range("L131").Value = 3396.6339
$ws.Range("N131").Value = -13476.6339

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 66359.03
$ws.Range("I132").Value = 43084.625
$ws.Range("K132").Value = 129253.875
$ws.Range("M132").Value = -126723.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 4500
$ws.Range("I5").Value = 4500
$ws.Range("K5").Value = 4500
$ws.Range("M5").Value = -4387

$ws.Range("H40").Value = 2469.85
$ws.Range("I40").Value = 2430.5386
$ws.Range("J40").Value = 2542.8572
$ws.Range("K40").Value = 2430.5386
$ws.Range("L40").Value = 2542.8572
$ws.Range("M40").Value = -2294.5386
$ws.Range("N40").Value = -2814.8572

$ws.Range("H93").Value = 316.2
$ws.Range("I93").Value = 320.25
$ws.Range("J93").Value = 300
$ws.Range("K93").Value = 320.25
$ws.Range("L93").Value = 300
$ws.Range("M93").Value = 927.75
$ws.Range("N93").Value = -2796

$ws.Range("H100").Value = 1687.75
$ws.Range("I100").Value = 1500.4
$ws.Range("K100").Value = 1500.4
$ws.Range("M100").Value = -959.4000000000001

$ws.Range("H122").Value = 3609.6843
$ws.Range("J122").Value = 3642.9412
$ws.Range("L122").Value = 10928.8236
$ws.Range("N122").Value = -15828.8236

$ws.Range("H132").Value = 24675
$ws.Range("I132").Value = 11021.875
$ws.Range("J132").Value = 94181.82000000001
$ws.Range("K132").Value = 33065.625
$ws.Range("L132").Value = 282545.46
$ws.Range("M132").Value = -30535.625
$ws.Range("N132").Value = -287605.46

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 19950
$ws.Range("J29").Value = 19950
$ws.Range("L29").Value = 19950
$ws.Range("M29").Value = -559.8

$ws.Range("H96").Value = 1890
$ws.Range("I96").Value = 1600.6
$ws.Range("J96").Value = 2179.4
$ws.Range("K96").Value = 1600.6
$ws.Range("L96").Value = 2179.4
$ws.Range("M96").Value = -227.5999999999999
$ws.Range("N96").Value = -4925.4

$ws.Range("H122").Value = 2452.8333
$ws.Range("I122").Value = 2128.7144
$ws.Range("J122").Value = 2906.6
$ws.Range("K122").Value = 6386.1432
$ws.Range("L122").Value = 8719.799999999999
$ws.Range("M122").Value = -3936.1432
$ws.Range("N122").Value = -13619.8

$ws.Range("H132").Value = 71195.5
$ws.Range("I132").Value = 59471.383
$ws.Range("J132").Value = 111057.5
$ws.Range("K132").Value = 178414.149
$ws.Range("L132").Value = 333172.5
$ws.Range("M132").Value = -175884.149
$ws.Range("N132").Value = -338232.5

$ws.Range("H136").Value = 49118.785
$ws.Range("I136").Value = 34167.766
$ws.Range("J136").Value = 86496.336
$ws.Range("K136").Value = 102503.298
$ws.Range("L136").Value = 259489.008
$ws.Range("M136").Value = -99953.29800000001
$ws.Range("N136").Value = -264589.008
